$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix export bug which bypasses deck image export:
# the "tagsEn" column (G) had blank/incorrect tag values for the robot
# cards and a stray "Construct" tag on the egg/dragon card.
$ws.Range("G5").Value = "Robot"
$ws.Range("G6").Value = "Robot"
$ws.Range("G7").Value = "Robot"
$ws.Range("G8").Value = "Robot"
$ws.Range("G9").Value = "Robot"
$ws.Range("G10").Value = "Dragon"

# Update the active window's view/selection to match where the author
# was working (row 9, column L).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L9").Select()
